# Update cryptocurrency price and volume(1h) data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.110.29"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.475.10"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.59"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.74"
$ws.Range("E6").Value = "  -1.85%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.066.05"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.476.14"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.101.58"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.41"
$ws.Range("E17").Value = "  -6.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.98"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.69"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.42"
$ws.Range("E20").Value = "  -2.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "385.51"
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.615.43"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.82"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.65"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.15"
$ws.Range("E30").Value = "  -4.45%  "
$ws.Range("E31").Value = "  -5.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("E32").Value = "  -4.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.505.32"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("E36").Value = "  -2.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.19"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("E39").Value = "  -3.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.13"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.29"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.82"
$ws.Range("E46").Value = "  -6.25%  "
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.12"
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.911"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.360.47"
$ws.Range("E51").Value = "  -4.67%  "
